$wb = $excel.ActiveWorkbook

# --- ExecData sheet: capitalize the action-type values in column C ---
$execSheet = $wb.Worksheets.Item("ExecData")
$execSheet.Range("C2").Value = "Navigate"
$execSheet.Range("C3").Value = "Input"
$execSheet.Range("C4").Value = "Click"

# Update selection on ExecData sheet and make it the active sheet/tab
$execSheet.Range("C5").Select() | Out-Null
$execSheet.Activate() | Out-Null
